$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 577, pushing existing rows 577-677 down to 579-679.
$ws.Rows("577:578").Insert()

# New row 577: weekly entry for Pepino ensalada (Primera), Región de Arica y Parinacota
$ws.Range("A577").Value = 6
$ws.Range("B577").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C577").Value = "Metropolitana"
$ws.Range("D577").Value = 45180
$ws.Range("E577").Value = 13
$ws.Range("F577").Value = 100112043
$ws.Range("G577").Value = "Pepino ensalada"
$ws.Range("H577").Value = "Sin especificar"
$ws.Range("I577").Value = "Primera"
$ws.Range("J577").Value = 1390
$ws.Range("K577").Value = 11000
$ws.Range("L577").Value = 12000
$ws.Range("M577").Value = 11482
$ws.Range("N577").Value = "`$/caja 60 unidades"
$ws.Range("O577").Value = "Región de Arica y Parinacota"
$ws.Range("P577").Value = 191
$ws.Range("Q577").Value = 60
$ws.Range("R577").Value = "Hortaliza"

# New row 578: weekly entry for Pepino ensalada (Segunda), Región de Arica y Parinacota
$ws.Range("A578").Value = 6
$ws.Range("B578").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C578").Value = "Metropolitana"
$ws.Range("D578").Value = 45180
$ws.Range("E578").Value = 13
$ws.Range("F578").Value = 100112043
$ws.Range("G578").Value = "Pepino ensalada"
$ws.Range("H578").Value = "Sin especificar"
$ws.Range("I578").Value = "Segunda"
$ws.Range("J578").Value = 220
$ws.Range("K578").Value = 10000
$ws.Range("L578").Value = 11000
$ws.Range("M578").Value = 10455
$ws.Range("N578").Value = "`$/caja 80 unidades"
$ws.Range("O578").Value = "Región de Arica y Parinacota"
$ws.Range("P578").Value = 131
$ws.Range("Q578").Value = 80
$ws.Range("R578").Value = "Hortaliza"
